$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value2 = 45129.50694444445
$ws.Range("B2").Value2 = 1.29
$ws.Range("C2").Value2 = 0.735
$ws.Range("D2").Value2 = 0.124
$ws.Range("E2").Value2 = 4.827
$ws.Range("F2").Value2 = 1.423
$ws.Range("G2").Value2 = 0.979
$ws.Range("H2").Value2 = 2.643
$ws.Range("I2").Value2 = 1.329
$ws.Range("J2").Value2 = 0.541
$ws.Range("K2").Value2 = 1.443
$ws.Range("L2").Value2 = 2.279
$ws.Range("M2").Value2 = 1.076
$ws.Range("N2").Value2 = 0.423
$ws.Range("O2").Value2 = 1.072
$ws.Range("P2").Value2 = 3.872
$ws.Range("Q2").Value2 = 0.514
$ws.Range("R2").Value2 = 0
$ws.Range("S2").Value2 = 0.001
$ws.Range("T2").Value2 = 14.386
$ws.Range("U2").Value2 = 3.712
$ws.Range("V2").Value2 = 2.121
$ws.Range("W2").Value2 = 1.552
$ws.Range("X2").Value2 = 1.645
$ws.Range("Y2").Value2 = 1.823
$ws.Range("Z2").Value2 = 2.314
$ws.Range("AA2").Value2 = 0.973
$ws.Range("AB2").Value2 = 0.721
$ws.Range("AC2").Value2 = 1.589
$ws.Range("AD2").Value2 = 1.435
$ws.Range("AE2").Value2 = 2.091
$ws.Range("AF2").Value2 = 2.771
$ws.Range("AG2").Value2 = 0.187
$ws.Range("AH2").Value2 = 2.45

# Row 3
$ws.Range("A3").Value2 = 45129.51388888889
$ws.Range("B3").Value2 = 7.839
$ws.Range("C3").Value2 = 5.792
$ws.Range("D3").Value2 = 0.346
$ws.Range("E3").Value2 = 17.923
$ws.Range("F3").Value2 = 13.714
$ws.Range("G3").Value2 = 6.272
$ws.Range("H3").Value2 = 19.203
$ws.Range("I3").Value2 = 9.483000000000001
$ws.Range("J3").Value2 = 4.259
$ws.Range("K3").Value2 = 6.671
$ws.Range("L3").Value2 = 7.493
$ws.Range("M3").Value2 = 7.249
$ws.Range("N3").Value2 = 2.042
$ws.Range("O3").Value2 = 6.128
$ws.Range("P3").Value2 = 9.875999999999999
$ws.Range("Q3").Value2 = 4.87
$ws.Range("R3").Value2 = 0.081
$ws.Range("S3").Value2 = 0.039
$ws.Range("T3").Value2 = 89.497
$ws.Range("U3").Value2 = 17.611
$ws.Range("V3").Value2 = 6.246
$ws.Range("W3").Value2 = 11.561
$ws.Range("X3").Value2 = 6.591
$ws.Range("Y3").Value2 = 1.568
$ws.Range("Z3").Value2 = 10.399
$ws.Range("AA3").Value2 = 5.159
$ws.Range("AB3").Value2 = 4.43
$ws.Range("AC3").Value2 = 5.627
$ws.Range("AD3").Value2 = 7.44
$ws.Range("AE3").Value2 = 0.8169999999999999
$ws.Range("AF3").Value2 = 17.387
$ws.Range("AG3").Value2 = 3.039
$ws.Range("AH3").Value2 = 7.722

# Row 4
$ws.Range("A4").Value2 = 45129.52083333334
$ws.Range("B4").Value2 = 3.143
$ws.Range("C4").Value2 = 2.3
$ws.Range("D4").Value2 = 0.183
$ws.Range("E4").Value2 = 7.385
$ws.Range("F4").Value2 = 5.401
$ws.Range("G4").Value2 = 2.573
$ws.Range("H4").Value2 = 12.021
$ws.Range("I4").Value2 = 3.819
$ws.Range("J4").Value2 = 1.793
$ws.Range("K4").Value2 = 2.761
$ws.Range("L4").Value2 = 3.209
$ws.Range("M4").Value2 = 2.89
$ws.Range("N4").Value2 = 0.851
$ws.Range("O4").Value2 = 2.443
$ws.Range("P4").Value2 = 4.322
$ws.Range("Q4").Value2 = 1.869
$ws.Range("R4").Value2 = 0.016
$ws.Range("S4").Value2 = 0
$ws.Range("T4").Value2 = 32.51
$ws.Range("U4").Value2 = 7.518
$ws.Range("V4").Value2 = 2.656
$ws.Range("W4").Value2 = 4.916
$ws.Range("X4").Value2 = 2.846
$ws.Range("Y4").Value2 = 0.8149999999999999
$ws.Range("Z4").Value2 = 6.163
$ws.Range("AA4").Value2 = 2.124
$ws.Range("AB4").Value2 = 1.785
$ws.Range("AC4").Value2 = 2.398
$ws.Range("AD4").Value2 = 3.065
$ws.Range("AE4").Value2 = 0.505
$ws.Range("AF4").Value2 = 11.494
$ws.Range("AG4").Value2 = 1.166
$ws.Range("AH4").Value2 = 3.267

# Row 5
$ws.Range("A5").Value2 = 45129.52777777778
$ws.Range("B5").Value2 = 11.36
$ws.Range("C5").Value2 = 8.49
$ws.Range("D5").Value2 = 0.45
$ws.Range("E5").Value2 = 25.12
$ws.Range("F5").Value2 = 20.25
$ws.Range("G5").Value2 = 9.029999999999999
$ws.Range("H5").Value2 = 31.88
$ws.Range("I5").Value2 = 13.78
$ws.Range("J5").Value2 = 6.21
$ws.Range("K5").Value2 = 9.390000000000001
$ws.Range("L5").Value2 = 10.29
$ws.Range("M5").Value2 = 10.49
$ws.Range("N5").Value2 = 2.89
$ws.Range("O5").Value2 = 8.880000000000001
$ws.Range("P5").Value2 = 13.24
$ws.Range("Q5").Value2 = 7.23
$ws.Range("R5").Value2 = 0.09
$ws.Range("S5").Value2 = 0.24
$ws.Range("T5").Value2 = 129.8
$ws.Range("U5").Value2 = 25.09
$ws.Range("V5").Value2 = 8.5
$ws.Range("W5").Value2 = 16.82
$ws.Range("X5").Value2 = 9.109999999999999
$ws.Range("Y5").Value2 = 1.53
$ws.Range("Z5").Value2 = 16.08
$ws.Range("AA5").Value2 = 7.34
$ws.Range("AB5").Value2 = 6.37
$ws.Range("AC5").Value2 = 7.74
$ws.Range("AD5").Value2 = 10.57
$ws.Range("AE5").Value2 = 0.37
$ws.Range("AF5").Value2 = 28.89
$ws.Range("AG5").Value2 = 4.59
$ws.Range("AH5").Value2 = 10.58

# Remove the old row 6 (no longer present in target range A1:AH5)
$ws.Rows.Item(6).Delete()

Write-Host "Done"